$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10, pushing the existing row 10
# (and its data) down to row 11. This preserves the old row's values/format
# while making room for the new weekly price entry.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly data.
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 44551
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100103
$ws.Cells.Item(10, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(10, 9).Value = 100103003
$ws.Cells.Item(10, 10).Value = "Damasco"
$ws.Cells.Item(10, 11).Value = "Castle Brite"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 120
$ws.Cells.Item(10, 14).Value = 15500
$ws.Cells.Item(10, 15).Value = 16000
$ws.Cells.Item(10, 16).Value = 15750
$ws.Cells.Item(10, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value = 1050
$ws.Cells.Item(10, 20).Value = 15
